$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text even when it looks like a number
# (e.g. "12.30", "1.00", "0.0000302") so Excel does not silently coerce it
# to a Double and lose the authored formatting. Restores "Normal" style
# afterwards so no stray NumberFormat/quotePrefix sticks to the cell.
function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).Style = "Normal"
}

$ws.Range("D2").Value = '69.297.22'
$ws.Range("E2").Value = '  -1.34%  '
$ws.Range("D3").Value = '3.506.55'
$ws.Range("E3").Value = '  -2.75%  '
Set-TextValue $ws "D5" '575.34'
$ws.Range("E5").Value = '  -1.06%  '
Set-TextValue $ws "D6" '185.54'
$ws.Range("E6").Value = '  -3.36%  '
$ws.Range("D7").Value = '3.496.97'
$ws.Range("E7").Value = '  -2.90%  '
Set-TextValue $ws "D8" '0.613'
$ws.Range("E8").Value = '  -3.27%  '
$ws.Range("E9").Value = '  +0.04%  '
$ws.Range("E10").Value = '  +4.37%  '
Set-TextValue $ws "D11" '0.647'
$ws.Range("E11").Value = '  -2.99%  '
Set-TextValue $ws "D12" '54.13'
$ws.Range("E12").Value = '  -3.41%  '
Set-TextValue $ws "D13" '0.0000302'
$ws.Range("E13").Value = '  -2.08%  '
Set-TextValue $ws "D14" '9.45'
$ws.Range("E14").Value = '  -2.86%  '
$ws.Range("D15").Value = '4.069.72'
$ws.Range("E15").Value = '  -2.82%  '
Set-TextValue $ws "D16" '19.36'
$ws.Range("E16").Value = '  -3.19%  '
$ws.Range("D17").Value = '69.262.38'
$ws.Range("E17").Value = '  -1.45%  '
$ws.Range("D18").Value = '3.502.48'
$ws.Range("E18").Value = '  -2.97%  '
Set-TextValue $ws "D19" '12.30'
$ws.Range("E19").Value = '  -3.41%  '
$ws.Range("E20").Value = '  -1.11%  '
Set-TextValue $ws "D21" '545.58'
$ws.Range("E21").Value = '  +12.73%  '
$ws.Range("E22").Value = '  -3.78%  '
Set-TextValue $ws "D23" '18.50'
$ws.Range("E23").Value = '  -3.65%  '
Set-TextValue $ws "D24" '4.97'
$ws.Range("E24").Value = '  -1.16%  '
Set-TextValue $ws "D25" '4.44'
$ws.Range("E25").Value = '  +0.64%  '
Set-TextValue $ws "D26" '94.07'
$ws.Range("E26").Value = '  -1.82%  '
Set-TextValue $ws "D27" '11.31'
$ws.Range("E27").Value = '  +1.58%  '
$ws.Range("E28").Value = '  -2.03%  '
Set-TextValue $ws "D29" '9.13'
$ws.Range("E29").Value = '  -2.98%  '
Set-TextValue $ws "D30" '31.85'
$ws.Range("E30").Value = '  -1.47%  '
Set-TextValue $ws "D31" '7.26'
$ws.Range("E31").Value = '  -6.70%  '
Set-TextValue $ws "D32" '12.60'
$ws.Range("E32").Value = '  +2.64%  '
Set-TextValue $ws "D33" '64.53'
$ws.Range("E33").Value = '  -3.44%  '
$ws.Range("E34").Value = '  -6.00%  '
Set-TextValue $ws "D35" '540.84'
$ws.Range("E35").Value = '  -8.46%  '
$ws.Range("B36").Value = 'Dai'
$ws.Range("C36").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextValue $ws "D36" '1.00'
$ws.Range("E36").Value = '  -0.05%  '
$ws.Range("B37").Value = 'Fetch.AI'
$ws.Range("C37").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws "D37" '3.08'
$ws.Range("E37").Value = '  +7.84%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws "D38" '37.95'
$ws.Range("E38").Value = '  -2.99%  '
$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue $ws "D39" '0.402'
$ws.Range("E39").Value = '  +0.83%  '
$ws.Range("E40").Value = '  -4.81%  '
Set-TextValue $ws "D41" '3.36'
$ws.Range("E41").Value = '  -3.08%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws "D42" '0.133'
$ws.Range("E42").Value = '  -2.75%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.309.52'
$ws.Range("E43").Value = '  +2.18%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue $ws "D44" '3.05'
$ws.Range("E44").Value = '  -8.00%  '
Set-TextValue $ws "D45" '2.98'
$ws.Range("E45").Value = '  -3.61%  '
Set-TextValue $ws "D46" '0.0445'
$ws.Range("E46").Value = '  -1.11%  '
$ws.Range("E47").Value = '  +3.91%  '
$ws.Range("E48").Value = '  -3.51%  '
Set-TextValue $ws "D49" '8.92'
$ws.Range("E49").Value = '  -6.56%  '
$ws.Range("E50").Value = '  -0.08%  '
Set-TextValue $ws "D51" '136.84'
$ws.Range("E51").Value = '  +1.93%  '
